# Applies the "1.4_APT122_FormativaFase1" edit:
#   - Section "2. Descripción breve del proyecto": rewrite body paragraph,
#     with a bold "pruebas de usabilidad" phrase inline.
#   - Section "5. Argumento sobre la factibilidad del proyecto": rewrite
#     body paragraph, with bold "BeautifulSoup", "proxies rotativos" and
#     "monitoreo continuo" phrases inline; append two new empty paragraphs
#     right after it.
#   - Section "6. Cumplimiento de los indicadores de calidad": rewrite
#     body paragraph, with a bold "actualización automática" phrase
#     inline; append one new empty paragraph right after it.

$d = $word.ActiveDocument

function Find-ParaIndexByText($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $needle) {
            return $i
        }
    }
    return -1
}

# Rewrites the text of paragraph $paraIndex using a list of
# @{text=...; bold=$true/$false} segments, concatenated in order, then
# re-applies bold formatting to the segments flagged bold=$true by
# searching for their text within the (now rewritten) paragraph.
function Set-ParaTextWithBoldSegments($paraIndex, $segments) {
    $full = ""
    foreach ($seg in $segments) { $full += $seg.text }
    $p = $d.Paragraphs.Item($paraIndex)
    $p.Range.Text = $full
    foreach ($seg in $segments) {
        if ($seg.bold) {
            $searchRange = $d.Paragraphs.Item($paraIndex).Range
            $searchRange.Find.ClearFormatting()
            $found = $searchRange.Find.Execute($seg.text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
            if ($found) {
                $searchRange.Bold = 1
            }
        }
    }
}

# Inserts a new, empty paragraph (matching the spacing of the paragraph it
# follows) right after paragraph $paraIndex. Returns the index of the
# newly created paragraph.
function Add-EmptyParagraphAfter($paraIndex) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $insertPoint = $r.End - 1
    $collapsed = $d.Range($insertPoint, $insertPoint)
    $collapsed.InsertParagraphAfter()
    return ($paraIndex + 1)
}

# ---------------------------------------------------------------------
# Process bottom-to-top so paragraph indices found earlier stay valid.
# ---------------------------------------------------------------------

# --- Section 6: "Cumplimiento de los indicadores de calidad" ---------
$segments6 = @(
    @{text="El proyecto cumple con los indicadores de calidad requeridos para el diseño del proyecto APT de acuerdo con los estándares de la disciplina. Cada fase del desarrollo, desde la arquitectura del sistema hasta las pruebas de seguridad, seguirá buenas prácticas industriales y académicas. También se ha planificado un sistema de "; bold=$false},
    @{text="actualización automática"; bold=$true},
    @{text=" para garantizar que los datos de la plataforma se mantengan actualizados y precisos en todo momento."; bold=$false}
)
$idx6 = Find-ParaIndexByText "6. Cumplimiento de los indicadores de calidad"
$body6 = $idx6 + 1
Set-ParaTextWithBoldSegments $body6 $segments6
Add-EmptyParagraphAfter $body6 | Out-Null

# --- Section 5: "Argumento sobre la factibilidad del proyecto" -------
$segments5 = @(
    @{text="El proyecto es completamente factible dentro del marco de la asignatura, ya que se desarrollará utilizando tecnologías accesibles y bien documentadas como Python para el web scraping, empleando una herramienta en específico como "; bold=$false},
    @{text="BeautifulSoup"; bold=$true},
    @{text=". Además, se implementarán estrategias para manejar posibles cambios en la estructura de los sitios web, como el uso de "; bold=$false},
    @{text="proxies rotativos"; bold=$true},
    @{text=" y "; bold=$false},
    @{text="monitoreo continuo"; bold=$true},
    @{text=" para evitar problemas a largo plazo relacionados con medidas anti-scraping. También, los protocolos estándar de ciberseguridad garantizarán la protección de los datos recopilados."; bold=$false}
)
$idx5 = Find-ParaIndexByText "5. Argumento sobre la factibilidad del proyecto"
$body5 = $idx5 + 1
Set-ParaTextWithBoldSegments $body5 $segments5
$afterBody5 = Add-EmptyParagraphAfter $body5
Add-EmptyParagraphAfter $afterBody5 | Out-Null

# --- Section 2: "Descripción breve del proyecto" ---------------------
$segments2 = @(
    @{text="La plataforma utiliza web scraping para obtener datos actualizados y relevantes sobre precios y características de repuestos, permitiendo a los usuarios ahorrar tiempo y realizar una búsqueda más precisa. Además, se ha decidido realizar "; bold=$false},
    @{text="pruebas de usabilidad"; bold=$true},
    @{text=" desde las fases tempranas del desarrollo para asegurar que la plataforma sea intuitiva y fácil de usar, alineándose con las expectativas de los usuarios."; bold=$false}
)
$idx2 = Find-ParaIndexByText "2. Descripción breve del proyecto"
$body2 = $idx2 + 1
Set-ParaTextWithBoldSegments $body2 $segments2

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
